$wb = $excel.ActiveWorkbook

# Mapping of row -> new F-column value for both "展览" and "全部类型" sheets
$updates = @{
    2  = 11722
    3  = 11349
    5  = 6
    6  = 1026
    11 = 10765
    12 = 4162
    13 = 16
    16 = 2467
    18 = 49
    20 = 449
    21 = 11142
    22 = 10925
    24 = 31
    27 = 28
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
